$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row to generic "question"/"answer" placeholders
$ws.Range("A1").Value = "question"
$ws.Range("B1").Value = "answer"

# Leave selection on B1, matching the final cursor position in the saved file
$ws.Range("B1").Select()
